$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 990
$ws.Range("J17").Value = 990
$ws.Range("L17").Value = 2970
$ws.Range("N17").Value = -3306
$ws.Range("H51").Value = 4334.8237
$ws.Range("I51").Value = 2186.25
$ws.Range("J51").Value = 6244.6665
$ws.Range("K51").Value = 2186.25
$ws.Range("L51").Value = 6244.6665
$ws.Range("M51").Value = -1702.25
$ws.Range("N51").Value = -7212.6665
$ws.Range("H64").Value = 4322.1763
$ws.Range("I64").Value = 2975
$ws.Range("J64").Value = 4501.8
$ws.Range("K64").Value = 2975
$ws.Range("L64").Value = 4501.8
$ws.Range("M64").Value = -2727
$ws.Range("N64").Value = -4997.8
$ws.Range("H67").Value = 4322.1763
$ws.Range("I67").Value = 2975
$ws.Range("J67").Value = 4501.8
$ws.Range("K67").Value = 2975
$ws.Range("L67").Value = 4501.8
$ws.Range("M67").Value = -2117
$ws.Range("N67").Value = -6217.8
$ws.Range("H98").Value = 1130.625
$ws.Range("I98").Value = 1130.625
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1130.625
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 367.375
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 1795.7
$ws.Range("I116").Value = 1694.625
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 1694.625
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1747.375
$ws.Range("N116").Value = -9084
$ws.Range("H122").Value = 1130.625
$ws.Range("I122").Value = 1130.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3391.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -941.875
$ws.Range("N122").ClearContents()
$ws.Range("H129").Value = 789.0333000000001
$ws.Range("I129").Value = 541
$ws.Range("J129").Value = 1037.0667
$ws.Range("K129").Value = 1623
$ws.Range("L129").Value = 3111.2001
$ws.Range("M129").Value = 3377
$ws.Range("N129").Value = -13111.2001
$ws.Range("H132").Value = 638891.4399999999
$ws.Range("I132").Value = 2274.879
$ws.Range("J132").Value = 4458591
$ws.Range("K132").Value = 6824.637
$ws.Range("L132").Value = 13375773
$ws.Range("M132").Value = -4294.637
$ws.Range("N132").Value = -13380833
$ws.Range("H137").Value = 1669289.8
$ws.Range("I137").Value = 2129615
$ws.Range("J137").Value = 5037.077
$ws.Range("K137").Value = 6388845
$ws.Range("L137").Value = 15111.231
$ws.Range("M137").Value = -6386295
$ws.Range("N137").Value = -20211.231

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 242.71428
$ws.Range("I5").Value = 233.16667
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 233.16667
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -121.16667
$ws.Range("N5").Value = -524
$ws.Range("H32").Value = 4348.84
$ws.Range("I32").Value = 4184.433
$ws.Range("J32").Value = 9664.666999999999
$ws.Range("K32").Value = 4184.433
$ws.Range("L32").Value = 9664.666999999999
$ws.Range("M32").Value = -3897.433
$ws.Range("N32").Value = -10238.667
$ws.Range("H61").Value = 37112484
$ws.Range("I61").Value = 41709680
$ws.Range("J61").Value = 334933.34
$ws.Range("K61").Value = 41709680
$ws.Range("L61").Value = 334933.34
$ws.Range("M61").Value = -41709468
$ws.Range("N61").Value = -335357.34
$ws.Range("H63").Value = 4012.125
$ws.Range("J63").Value = 3299.5
$ws.Range("L63").Value = 3299.5
$ws.Range("N63").Value = -4671.5
$ws.Range("H66").Value = 4012.125
$ws.Range("J66").Value = 3299.5
$ws.Range("L66").Value = 16497.5
$ws.Range("N66").Value = -23361.5
$ws.Range("H132").Value = 7174248.5
$ws.Range("I132").Value = 8214406.5
$ws.Range("J132").Value = 124288
$ws.Range("K132").Value = 24643219.5
$ws.Range("L132").Value = 372864
$ws.Range("M132").Value = -24640689.5
$ws.Range("N132").Value = -377924
$ws.Range("H136").Value = 37112484
$ws.Range("I136").Value = 41709680
$ws.Range("J136").Value = 334933.34
$ws.Range("K136").Value = 125129040
$ws.Range("L136").Value = 1004800.02
$ws.Range("M136").Value = -125126490
$ws.Range("N136").Value = -1009900.02

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 242.71428
$ws.Range("I4").Value = 233.16667
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 233.16667
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -118.16667
$ws.Range("N4").Value = -530
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H82").Value = 9115.75
$ws.Range("I82").Value = 4882.6
$ws.Range("J82").Value = 30281.5
$ws.Range("K82").Value = 4882.6
$ws.Range("L82").Value = 30281.5
$ws.Range("M82").Value = -4499.6
$ws.Range("N82").Value = -31047.5
$ws.Range("H85").Value = 9115.75
$ws.Range("I85").Value = 4882.6
$ws.Range("J85").Value = 30281.5
$ws.Range("K85").Value = 4882.6
$ws.Range("L85").Value = 30281.5
$ws.Range("M85").Value = -3556.6
$ws.Range("N85").Value = -32933.5
$ws.Range("H107").Value = 4652.96
$ws.Range("I107").Value = 3980.2
$ws.Range("J107").Value = 5662.1
$ws.Range("K107").Value = 3980.2
$ws.Range("L107").Value = 5662.1
$ws.Range("M107").Value = -2060.2
$ws.Range("N107").Value = -9502.1
$ws.Range("H134").Value = 8335701
$ws.Range("I134").Value = 2026.5172
$ws.Range("K134").Value = 6079.5516
$ws.Range("M134").Value = -3544.5516

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44119.38
$ws.Range("I31").Value = 32176.395
$ws.Range("J31").Value = 67302.82000000001
$ws.Range("K31").Value = 32176.395
$ws.Range("L31").Value = 67302.82000000001
$ws.Range("M31").Value = -31881.395
$ws.Range("N31").Value = -67892.82000000001
$ws.Range("H34").Value = 44119.38
$ws.Range("I34").Value = 32176.395
$ws.Range("J34").Value = 67302.82000000001
$ws.Range("K34").Value = 32176.395
$ws.Range("L34").Value = 67302.82000000001
$ws.Range("M34").Value = -31974.395
$ws.Range("N34").Value = -67706.82000000001
$ws.Range("H132").Value = 24920.303
$ws.Range("I132").Value = 1428.6364
$ws.Range("K132").Value = 4285.9092
$ws.Range("M132").Value = -1755.9092
$ws.Range("H134").Value = 27841.875
$ws.Range("I134").Value = 759.6070999999999
$ws.Range("J134").Value = 91033.836
$ws.Range("K134").Value = 2278.8213
$ws.Range("L134").Value = 273101.508
$ws.Range("M134").Value = 256.1787000000004
$ws.Range("N134").Value = -278171.508

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 822.39703
$ws.Range("I131").Value = 348.30768
$ws.Range("J131").Value = 934.4545000000001
$ws.Range("K131").Value = 1044.92304
$ws.Range("L131").Value = 2803.3635
$ws.Range("M131").Value = 3995.07696
$ws.Range("N131").Value = -12883.3635
$ws.Range("H140").Value = 2776.195
$ws.Range("I140").Value = 3230.8333
$ws.Range("J140").Value = 2134.353
$ws.Range("K140").Value = 9692.499899999999
$ws.Range("L140").Value = 6403.059
$ws.Range("M140").Value = -4512.499899999999
$ws.Range("N140").Value = -16763.059
$ws.Range("H141").Value = 8855.625
$ws.Range("I141").Value = 3790
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 11370
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -6190
$ws.Range("N141").Value = -70360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25844.467
$ws.Range("I132").Value = 2953.2856
$ws.Range("J132").Value = 63547.59
$ws.Range("K132").Value = 8859.856800000001
$ws.Range("L132").Value = 190642.77
$ws.Range("M132").Value = -6329.856800000001
$ws.Range("N132").Value = -195702.77
$ws.Range("H136").Value = 31435.857
$ws.Range("I136").Value = 18666.826
$ws.Range("J136").Value = 146357.14
$ws.Range("K136").Value = 56000.478
$ws.Range("L136").Value = 439071.42
$ws.Range("M136").Value = -53450.478
$ws.Range("N136").Value = -444171.42

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1573.5
$ws.Range("I126").Value = 681.3333
$ws.Range("J126").Value = 4250
$ws.Range("K126").Value = 2043.9999
$ws.Range("L126").Value = 12750
$ws.Range("M126").Value = 426.0001
$ws.Range("N126").Value = -17690
$ws.Range("H132").Value = 40051.71
$ws.Range("I132").Value = 35591.07
$ws.Range("J132").Value = 45676
$ws.Range("K132").Value = 106773.21
$ws.Range("L132").Value = 137028
$ws.Range("M132").Value = -104243.21
$ws.Range("N132").Value = -142088
$ws.Range("H136").Value = 31644.238
$ws.Range("I136").Value = 20220.27
$ws.Range("J136").Value = 71247.336
$ws.Range("K136").Value = 60660.81
$ws.Range("L136").Value = 213742.008
$ws.Range("M136").Value = -58110.81
$ws.Range("N136").Value = -218842.008
